$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D3").Value = -0.0313
$ws.Range("E2:E3").Value = 0.0673
$ws.Range("G2:G3").Value = 1.088757396449704
$ws.Range("H2:H3").Value = 1.088757396449704
$ws.Range("I2:I3").Value = 0.5798816568047338
$ws.Range("J2:J3").Value = 0.5149842628729701
$ws.Range("K2:K3").Value = 20.9
$ws.Range("L2:L3").Value = 0.6183431952662722
$ws.Range("M2:M3").Value = 8.33
$ws.Range("N2:N3").Value = 0.02679318108716629
$ws.Range("O2:O3").Value = 0.3985645933014355
$ws.Range("P2:P3").Value = 8.33
$ws.Range("Q2:Q3").Value = 0.02679318108716629
$ws.Range("R2:R3").Value = 0.3985645933014355
$ws.Range("U2:U3").Value = 38.5
$ws.Range("V2:V3").Value = 0.1238340302348022
$ws.Range("W2:W3").Value = 0.06959706959706959
$ws.Range("X2:X3").Value = 0.04031256956053166
$ws.Range("Y2:Y3").Value = 0.02928450003653793
$ws.Range("Z2:Z3").Value = 0.1283511809827599
$ws.Range("AA2:AA3").Value = 0.06609883832728178
$ws.Range("AB2:AB3").Value = 0.04020795579160034
$ws.Range("AC2:AC3").Value = 0.02589088253568143
$ws.Range("AD2:AD3").Value = 1.34
$ws.Range("AE2:AE3").Value = 0
$ws.Range("AF2:AF3").Value = 1.34
$ws.Range("AG2:AG3").Value = -37.16
$ws.Range("AH2:AH3").Value = 0.004291570586728159
$ws.Range("AI2:AI3").Value = 0.00458213650663384
$ws.Range("AJ2:AJ3").Value = -0.1357492511141959
$ws.Range("AK2:AK3").Value = -0.1463337796329841
$ws.Range("AN2:AN3").Value = 0.06767676767676768
$ws.Range("AP2:AP3").Value = -1.876767676767676
